$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 339, shifting existing rows 339:497 down to 340:498.
$ws.Rows(339).Insert()

# Populate the newly inserted row 339 with the new weekly data entry.
$ws.Range("A339").Value = 9
$ws.Range("B339").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C339").Value = "Metropolitana"
$ws.Range("D339").Value = 45205
$ws.Range("E339").Value = 13
$ws.Range("F339").Value = 300000001
$ws.Range("G339").Value = "Rabanito"
$ws.Range("H339").Value = "Sin especificar"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 7000
$ws.Range("K339").Value = 3000
$ws.Range("L339").Value = 3000
$ws.Range("M339").Value = 3000
$ws.Range("N339").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O339").Value = "Provincia de Chacabuco"
$ws.Range("P339").Value = 30
$ws.Range("Q339").Value = 100
$ws.Range("R339").Value = "Hortaliza"

# Match the date cell style used by the other date cells in column D.
$ws.Range("D339").NumberFormat = $ws.Range("D340").NumberFormat
